$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.879.20"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "3.070.41"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.40"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.01"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").Value = "3.596.89"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "57.867.61"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "3.067.15"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.14"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.06"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.37"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.500"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.01"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "0.0₃0905"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.76"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.70"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.55"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.53"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0674"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "3.111.80"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.53"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.655"
$ws.Range("D44").Value = "2.277.98"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.59"
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.92"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.936"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  +6.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "255.50"
$ws.Range("E51").Value = "  +8.27%  "
